$wb = $excel.ActiveWorkbook

# "Generate Report for handoff" - update the Latest Handoff Datetime for the
# 9e2f9f8b-6d06-4fdc-b8db-dab61782ac03 entry (row 5) on each language sheet,
# reflecting a newly generated handoff report/timestamp.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-15 14:03:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-15 14:03:29"
